$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "RM 232" row (row 26) and the "SC 92" row (originally row 28,
# but after the first deletion it shifts up to row 27).
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# Fill in the previously-missing F value for "SC 120" (now row 30)
$ws.Range("F30").Value = 16.89

# Clear the F value for "SC 193" (now row 32) back to blank/missing
$ws.Range("F32").Value = ""
